$d = $word.ActiveDocument

# --- Part 1: collapse the "block width" sentence runs and drop the stray
#     _GoBack bookmark that Word had left around them -------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$r1 = $d.Content
[void]$r1.Find.Execute("The best value for block width was always 64. The lower value was ", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "The best value for block width was always 64. The lower value was ", 2)

$r2 = $d.Content
[void]$r2.Find.Execute("32,  the", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "32,  the", 2)

$r3 = $d.Content
[void]$r3.Find.Execute(" higher value was 128. ", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    " higher value was 128. ", 2)

# --- Part 2: append a new numbered list paragraph with the repo link,
#     carrying the _GoBack bookmark with it ----------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "https://github.com/JHobbie/cs392hw4test/commits/master"

# Find an existing numbered-list template in the document and apply it
# to the new paragraph so it matches numId=1 / ilvl=0.
$listTemplate = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.ListFormat.ListType -eq 3) {
        $listTemplate = $p.Range.ListFormat.ListTemplate
        break
    }
}
$newPara.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 2, $false, 0)

$bmRange = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
